$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 2-5: descriptive/source text (not date-like, plain assignment is safe)
$ws.Range("A2").Value = '                  Central Bank of Sri Lanka'
$ws.Range("A3").Value = '                  National Gem and Jewellery Authority'
$ws.Range("A4").Value = '                  Sri Lanka Customs'
$ws.Range("A5").Value = '  Sources: Ceylon Petroleum Corporation and Other Exporters of Petroleum'

# Rows 6-10: 'YYYY Month' strings read as dates by Excel's type-inference,
# so force the cell format to Text first to keep the literal string.
$ws.Range("A6:A10").NumberFormat = "@"
$ws.Range("A6").Value = '2006 January'
$ws.Range("A7").Value = '2007 January'
$ws.Range("A8").Value = '2008 January'
$ws.Range("A9").Value = '2009 January'
$ws.Range("A10").Value = '2010 January'

# Rows 11-70: month names / labels / table titles (plain text, safe to assign directly)
$ws.Range("A11").Value = 'April'
$ws.Range("A12").Value = 'April'
$ws.Range("A13").Value = 'April'
$ws.Range("A14").Value = 'April'
$ws.Range("A15").Value = 'April'
$ws.Range("A16").Value = 'August'
$ws.Range("A17").Value = 'August'
$ws.Range("A18").Value = 'August'
$ws.Range("A19").Value = 'August'
$ws.Range("A20").Value = 'August'
$ws.Range("A21").Value = 'December'
$ws.Range("A22").Value = 'December'
$ws.Range("A23").Value = 'December'
$ws.Range("A24").Value = 'December'
$ws.Range("A25").Value = 'December'
$ws.Range("A26").Value = 'February'
$ws.Range("A27").Value = 'February'
$ws.Range("A28").Value = 'February'
$ws.Range("A29").Value = 'February'
$ws.Range("A30").Value = 'February'
$ws.Range("A31").Value = 'Industrial Exports'
$ws.Range("A32").Value = 'July'
$ws.Range("A33").Value = 'July'
$ws.Range("A34").Value = 'July'
$ws.Range("A35").Value = 'July'
$ws.Range("A36").Value = 'July'
$ws.Range("A37").Value = 'June'
$ws.Range("A38").Value = 'June'
$ws.Range("A39").Value = 'June'
$ws.Range("A40").Value = 'June'
$ws.Range("A41").Value = 'June'
$ws.Range("A42").Value = 'March'
$ws.Range("A43").Value = 'March'
$ws.Range("A44").Value = 'March'
$ws.Range("A45").Value = 'March'
$ws.Range("A46").Value = 'March'
$ws.Range("A47").Value = 'May'
$ws.Range("A48").Value = 'May'
$ws.Range("A49").Value = 'May'
$ws.Range("A50").Value = 'May'
$ws.Range("A51").Value = 'May'
$ws.Range("A52").Value = 'November'
$ws.Range("A53").Value = 'November'
$ws.Range("A54").Value = 'November'
$ws.Range("A55").Value = 'November'
$ws.Range("A56").Value = 'November'
$ws.Range("A57").Value = 'October'
$ws.Range("A58").Value = 'October'
$ws.Range("A59").Value = 'October'
$ws.Range("A60").Value = 'October'
$ws.Range("A61").Value = 'October'
$ws.Range("A62").Value = 'Period'
$ws.Range("A63").Value = 'September'
$ws.Range("A64").Value = 'September'
$ws.Range("A65").Value = 'September'
$ws.Range("A66").Value = 'September'
$ws.Range("A67").Value = 'September'
$ws.Range("A68").Value = 'Table 2.02.6: Exports (Rupees Million)'
$ws.Range("A69").Value = 'Table 2.02: Exports - Monthly (2006-2010)'
$ws.Range("A70").Value = 'Table 2.02: Exports - Monthly (2006-2010)'

# Rows 71-223 previously held the leftover sorted/duplicated date & label values;
# clear them so the column ends at row 70 (rows 224+ were already empty).
$ws.Range("A71:A223").ClearContents()

